$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right above row 341 - this shifts the existing data
# (rows 341:430) down to (344:433), matching the new dimension A1:T433.
$ws.Rows("341:343").Insert()

# Populate the newly inserted rows with the new weekly price block
# (Fecha = 2021-12-13 / serial 44543), keeping the same row layout as
# every other block: Pintón / Primera Maduro / Primera Pintón.

# Row 341 - Pintón
$ws.Range("A341").Value = 8
$ws.Range("B341").Value = "Terminal La Palmera de La Serena"
$ws.Range("C341").Value = "Coquimbo"
$ws.Range("D341").Value = 44543
$ws.Range("E341").Value = 4
$ws.Range("F341").Value = "Fruta"
$ws.Range("G341").Value = 100108
$ws.Range("H341").Value = "Tropicales y subtropicales"
$ws.Range("I341").Value = 100108006
$ws.Range("J341").Value = "Plátano"
$ws.Range("K341").Value = "Sin especificar"
$ws.Range("L341").Value = "Pintón"
$ws.Range("M341").Value = 80
$ws.Range("N341").Value = 14000
$ws.Range("O341").Value = 14000
$ws.Range("P341").Value = 14000
$ws.Range("Q341").Value = "$/caja 20 kilos"
$ws.Range("R341").Value = "Ecuador"
$ws.Range("S341").Value = 700
$ws.Range("T341").Value = 20

# Row 342 - Primera Maduro
$ws.Range("A342").Value = 8
$ws.Range("B342").Value = "Terminal La Palmera de La Serena"
$ws.Range("C342").Value = "Coquimbo"
$ws.Range("D342").Value = 44543
$ws.Range("E342").Value = 4
$ws.Range("F342").Value = "Fruta"
$ws.Range("G342").Value = 100108
$ws.Range("H342").Value = "Tropicales y subtropicales"
$ws.Range("I342").Value = 100108006
$ws.Range("J342").Value = "Plátano"
$ws.Range("K342").Value = "Sin especificar"
$ws.Range("L342").Value = "Primera Maduro"
$ws.Range("M342").Value = 120
$ws.Range("N342").Value = 16000
$ws.Range("O342").Value = 16000
$ws.Range("P342").Value = 16000
$ws.Range("Q342").Value = "$/caja 20 kilos"
$ws.Range("R342").Value = "Ecuador"
$ws.Range("S342").Value = 800
$ws.Range("T342").Value = 20

# Row 343 - Primera Pintón
$ws.Range("A343").Value = 8
$ws.Range("B343").Value = "Terminal La Palmera de La Serena"
$ws.Range("C343").Value = "Coquimbo"
$ws.Range("D343").Value = 44543
$ws.Range("E343").Value = 4
$ws.Range("F343").Value = "Fruta"
$ws.Range("G343").Value = 100108
$ws.Range("H343").Value = "Tropicales y subtropicales"
$ws.Range("I343").Value = 100108006
$ws.Range("J343").Value = "Plátano"
$ws.Range("K343").Value = "Sin especificar"
$ws.Range("L343").Value = "Primera Pintón"
$ws.Range("M343").Value = 120
$ws.Range("N343").Value = 17000
$ws.Range("O343").Value = 17000
$ws.Range("P343").Value = 17000
$ws.Range("Q343").Value = "$/caja 20 kilos"
$ws.Range("R343").Value = "Ecuador"
$ws.Range("S343").Value = 850
$ws.Range("T343").Value = 20

# Ensure the date column keeps the same custom date style/number format
# used throughout column D (style index 2 in the original file).
$ws.Range("D341:D343").NumberFormat = "YYYY-MM-DD HH:MM:SS"
